$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data scraped from coinranking.com.
# Force text NumberFormat before assigning so that numeric-looking values
# (e.g. "0.9990", "26.829.50") are preserved exactly as text, matching the
# original inline-string cell contents, then clear the format override so
# the cell's style index returns to the sheet's default (unstyled) state.
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "26.829.50"
$ws.Cells.Item(2, 4).ClearFormats()
$ws.Cells.Item(2, 5).NumberFormat = "@"
$ws.Cells.Item(2, 5).Value = "  -2.20%  "
$ws.Cells.Item(2, 5).ClearFormats()
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "1.800.93"
$ws.Cells.Item(3, 4).ClearFormats()
$ws.Cells.Item(3, 5).NumberFormat = "@"
$ws.Cells.Item(3, 5).Value = "  -1.47%  "
$ws.Cells.Item(3, 5).ClearFormats()
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "0.9990"
$ws.Cells.Item(4, 4).ClearFormats()
$ws.Cells.Item(4, 5).NumberFormat = "@"
$ws.Cells.Item(4, 5).Value = "  -0.12%  "
$ws.Cells.Item(4, 5).ClearFormats()
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "308.90"
$ws.Cells.Item(5, 4).ClearFormats()
$ws.Cells.Item(5, 5).NumberFormat = "@"
$ws.Cells.Item(5, 5).Value = "  -1.70%  "
$ws.Cells.Item(5, 5).ClearFormats()
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "0.9993"
$ws.Cells.Item(6, 4).ClearFormats()
$ws.Cells.Item(6, 5).NumberFormat = "@"
$ws.Cells.Item(6, 5).Value = "  -0.07%  "
$ws.Cells.Item(6, 5).ClearFormats()
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.4635"
$ws.Cells.Item(7, 4).ClearFormats()
$ws.Cells.Item(7, 5).NumberFormat = "@"
$ws.Cells.Item(7, 5).Value = "  +3.65%  "
$ws.Cells.Item(7, 5).ClearFormats()
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.3731"
$ws.Cells.Item(8, 4).ClearFormats()
$ws.Cells.Item(8, 5).NumberFormat = "@"
$ws.Cells.Item(8, 5).Value = "  -0.76%  "
$ws.Cells.Item(8, 5).ClearFormats()
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.07286"
$ws.Cells.Item(9, 4).ClearFormats()
$ws.Cells.Item(9, 5).NumberFormat = "@"
$ws.Cells.Item(9, 5).Value = "  -3.20%  "
$ws.Cells.Item(9, 5).ClearFormats()
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.8595"
$ws.Cells.Item(10, 4).ClearFormats()
$ws.Cells.Item(10, 5).NumberFormat = "@"
$ws.Cells.Item(10, 5).Value = "  -3.60%  "
$ws.Cells.Item(10, 5).ClearFormats()
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "20.42"
$ws.Cells.Item(11, 4).ClearFormats()
$ws.Cells.Item(11, 5).NumberFormat = "@"
$ws.Cells.Item(11, 5).Value = "  -3.09%  "
$ws.Cells.Item(11, 5).ClearFormats()
$ws.Cells.Item(12, 2).Value = "WrappedEther"
$ws.Cells.Item(12, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "1.746.18"
$ws.Cells.Item(12, 4).ClearFormats()
$ws.Cells.Item(12, 5).NumberFormat = "@"
$ws.Cells.Item(12, 5).Value = "  -4.15%  "
$ws.Cells.Item(12, 5).ClearFormats()
$ws.Cells.Item(13, 2).Value = "Polkadot"
$ws.Cells.Item(13, 3).Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "5.327"
$ws.Cells.Item(13, 4).ClearFormats()
$ws.Cells.Item(13, 5).NumberFormat = "@"
$ws.Cells.Item(13, 5).Value = "  -1.57%  "
$ws.Cells.Item(13, 5).ClearFormats()
$ws.Cells.Item(14, 2).Value = "Chainlink"
$ws.Cells.Item(14, 3).Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "6.497"
$ws.Cells.Item(14, 4).ClearFormats()
$ws.Cells.Item(14, 5).NumberFormat = "@"
$ws.Cells.Item(14, 5).Value = "  -3.74%  "
$ws.Cells.Item(14, 5).ClearFormats()
$ws.Cells.Item(15, 2).Value = "TRON"
$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "0.07044"
$ws.Cells.Item(15, 4).ClearFormats()
$ws.Cells.Item(15, 5).NumberFormat = "@"
$ws.Cells.Item(15, 5).Value = "  -1.12%  "
$ws.Cells.Item(15, 5).ClearFormats()
$ws.Cells.Item(16, 2).Value = "Litecoin"
$ws.Cells.Item(16, 3).Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "91.27"
$ws.Cells.Item(16, 4).ClearFormats()
$ws.Cells.Item(16, 5).NumberFormat = "@"
$ws.Cells.Item(16, 5).Value = "  -3.30%  "
$ws.Cells.Item(16, 5).ClearFormats()
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "0.9996"
$ws.Cells.Item(17, 4).ClearFormats()
$ws.Cells.Item(17, 5).NumberFormat = "@"
$ws.Cells.Item(17, 5).Value = "  -0.15%  "
$ws.Cells.Item(17, 5).ClearFormats()
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "0.000008663"
$ws.Cells.Item(18, 4).ClearFormats()
$ws.Cells.Item(18, 5).NumberFormat = "@"
$ws.Cells.Item(18, 5).Value = "  -1.82%  "
$ws.Cells.Item(18, 5).ClearFormats()
$ws.Cells.Item(19, 5).NumberFormat = "@"
$ws.Cells.Item(19, 5).Value = "  +0.00%  "
$ws.Cells.Item(19, 5).ClearFormats()
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "14.67"
$ws.Cells.Item(20, 4).ClearFormats()
$ws.Cells.Item(20, 5).NumberFormat = "@"
$ws.Cells.Item(20, 5).Value = "  -3.72%  "
$ws.Cells.Item(20, 5).ClearFormats()
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "26.802.05"
$ws.Cells.Item(21, 4).ClearFormats()
$ws.Cells.Item(21, 5).NumberFormat = "@"
$ws.Cells.Item(21, 5).Value = "  -2.30%  "
$ws.Cells.Item(21, 5).ClearFormats()
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "5.299"
$ws.Cells.Item(22, 4).ClearFormats()
$ws.Cells.Item(22, 5).NumberFormat = "@"
$ws.Cells.Item(22, 5).Value = "  +0.56%  "
$ws.Cells.Item(22, 5).ClearFormats()
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "10.67"
$ws.Cells.Item(23, 4).ClearFormats()
$ws.Cells.Item(23, 5).NumberFormat = "@"
$ws.Cells.Item(23, 5).Value = "  -2.18%  "
$ws.Cells.Item(23, 5).ClearFormats()
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "1.916.43"
$ws.Cells.Item(24, 4).ClearFormats()
$ws.Cells.Item(24, 5).NumberFormat = "@"
$ws.Cells.Item(24, 5).Value = "  -6.70%  "
$ws.Cells.Item(24, 5).ClearFormats()
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "1.904"
$ws.Cells.Item(25, 4).ClearFormats()
$ws.Cells.Item(25, 5).NumberFormat = "@"
$ws.Cells.Item(25, 5).Value = "  -3.72%  "
$ws.Cells.Item(25, 5).ClearFormats()
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "150.56"
$ws.Cells.Item(26, 4).ClearFormats()
$ws.Cells.Item(26, 5).NumberFormat = "@"
$ws.Cells.Item(26, 5).Value = "  -0.56%  "
$ws.Cells.Item(26, 5).ClearFormats()
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "18.22"
$ws.Cells.Item(27, 4).ClearFormats()
$ws.Cells.Item(27, 5).NumberFormat = "@"
$ws.Cells.Item(27, 5).Value = "  -2.17%  "
$ws.Cells.Item(27, 5).ClearFormats()
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "2.141"
$ws.Cells.Item(28, 4).ClearFormats()
$ws.Cells.Item(28, 5).NumberFormat = "@"
$ws.Cells.Item(28, 5).Value = "  -9.79%  "
$ws.Cells.Item(28, 5).ClearFormats()
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "5.240"
$ws.Cells.Item(29, 4).ClearFormats()
$ws.Cells.Item(29, 5).NumberFormat = "@"
$ws.Cells.Item(29, 5).Value = "  -2.49%  "
$ws.Cells.Item(29, 5).ClearFormats()
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "114.29"
$ws.Cells.Item(30, 4).ClearFormats()
$ws.Cells.Item(30, 5).NumberFormat = "@"
$ws.Cells.Item(30, 5).Value = "  -2.86%  "
$ws.Cells.Item(30, 5).ClearFormats()
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "0.08901"
$ws.Cells.Item(31, 4).ClearFormats()
$ws.Cells.Item(31, 5).NumberFormat = "@"
$ws.Cells.Item(31, 5).Value = "  +0.68%  "
$ws.Cells.Item(31, 5).ClearFormats()
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "0.7621"
$ws.Cells.Item(32, 4).ClearFormats()
$ws.Cells.Item(32, 5).NumberFormat = "@"
$ws.Cells.Item(32, 5).Value = "  -2.65%  "
$ws.Cells.Item(32, 5).ClearFormats()
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "1.167"
$ws.Cells.Item(33, 4).ClearFormats()
$ws.Cells.Item(33, 5).NumberFormat = "@"
$ws.Cells.Item(33, 5).Value = "  -3.04%  "
$ws.Cells.Item(33, 5).ClearFormats()
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "4.445"
$ws.Cells.Item(34, 4).ClearFormats()
$ws.Cells.Item(34, 5).NumberFormat = "@"
$ws.Cells.Item(34, 5).Value = "  -1.47%  "
$ws.Cells.Item(34, 5).ClearFormats()
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "2.889"
$ws.Cells.Item(35, 4).ClearFormats()
$ws.Cells.Item(35, 5).NumberFormat = "@"
$ws.Cells.Item(35, 5).Value = "  +0.04%  "
$ws.Cells.Item(35, 5).ClearFormats()
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "0.9992"
$ws.Cells.Item(36, 4).ClearFormats()
$ws.Cells.Item(36, 5).NumberFormat = "@"
$ws.Cells.Item(36, 5).Value = "  -0.04%  "
$ws.Cells.Item(36, 5).ClearFormats()
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "1.118"
$ws.Cells.Item(37, 4).ClearFormats()
$ws.Cells.Item(37, 5).NumberFormat = "@"
$ws.Cells.Item(37, 5).Value = "  +0.76%  "
$ws.Cells.Item(37, 5).ClearFormats()
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.01946"
$ws.Cells.Item(38, 4).ClearFormats()
$ws.Cells.Item(38, 5).NumberFormat = "@"
$ws.Cells.Item(38, 5).Value = "  -2.26%  "
$ws.Cells.Item(38, 5).ClearFormats()
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.05215"
$ws.Cells.Item(39, 4).ClearFormats()
$ws.Cells.Item(39, 5).NumberFormat = "@"
$ws.Cells.Item(39, 5).Value = "  -2.04%  "
$ws.Cells.Item(39, 5).ClearFormats()
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "2.387"
$ws.Cells.Item(40, 4).ClearFormats()
$ws.Cells.Item(40, 5).NumberFormat = "@"
$ws.Cells.Item(40, 5).Value = "  +4.66%  "
$ws.Cells.Item(40, 5).ClearFormats()
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "2.901"
$ws.Cells.Item(41, 4).ClearFormats()
$ws.Cells.Item(41, 5).NumberFormat = "@"
$ws.Cells.Item(41, 5).Value = "  +1.28%  "
$ws.Cells.Item(41, 5).ClearFormats()
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "7.181"
$ws.Cells.Item(42, 4).ClearFormats()
$ws.Cells.Item(42, 5).NumberFormat = "@"
$ws.Cells.Item(42, 5).Value = "  -2.64%  "
$ws.Cells.Item(42, 5).ClearFormats()
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.5221"
$ws.Cells.Item(43, 4).ClearFormats()
$ws.Cells.Item(43, 5).NumberFormat = "@"
$ws.Cells.Item(43, 5).Value = "  -1.68%  "
$ws.Cells.Item(43, 5).ClearFormats()
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.1654"
$ws.Cells.Item(44, 4).ClearFormats()
$ws.Cells.Item(44, 5).NumberFormat = "@"
$ws.Cells.Item(44, 5).Value = "  -4.29%  "
$ws.Cells.Item(44, 5).ClearFormats()
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "8.524"
$ws.Cells.Item(45, 4).ClearFormats()
$ws.Cells.Item(45, 5).NumberFormat = "@"
$ws.Cells.Item(45, 5).Value = "  -2.91%  "
$ws.Cells.Item(45, 5).ClearFormats()
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "0.5025"
$ws.Cells.Item(46, 4).ClearFormats()
$ws.Cells.Item(46, 5).NumberFormat = "@"
$ws.Cells.Item(46, 5).Value = "  -2.05%  "
$ws.Cells.Item(46, 5).ClearFormats()
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "10.21"
$ws.Cells.Item(47, 4).ClearFormats()
$ws.Cells.Item(47, 5).NumberFormat = "@"
$ws.Cells.Item(47, 5).Value = "  -4.91%  "
$ws.Cells.Item(47, 5).ClearFormats()
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "104.32"
$ws.Cells.Item(48, 4).ClearFormats()
$ws.Cells.Item(48, 5).NumberFormat = "@"
$ws.Cells.Item(48, 5).Value = "  -1.42%  "
$ws.Cells.Item(48, 5).ClearFormats()
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "0.9989"
$ws.Cells.Item(49, 4).ClearFormats()
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "1.651"
$ws.Cells.Item(50, 4).ClearFormats()
$ws.Cells.Item(50, 5).NumberFormat = "@"
$ws.Cells.Item(50, 5).Value = "  -3.16%  "
$ws.Cells.Item(50, 5).ClearFormats()
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "0.06286"
$ws.Cells.Item(51, 4).ClearFormats()
$ws.Cells.Item(51, 5).NumberFormat = "@"
$ws.Cells.Item(51, 5).Value = "  -1.42%  "
$ws.Cells.Item(51, 5).ClearFormats()
